$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "ab2cd12983qwe@gmail.com"
$ws.Range("D3").Value = "lo982ki87hgf3@gmail.com"
$ws.Range("D4").Value = "lok2793ijuh76tg@hotmail.com"
$ws.Range("D5").Value = "ada17mg124@gmail.com"
$ws.Range("D6").Value = "as23dq678wer@yahoo.com"

$ws.Range("D6").Select()
